# Replace the 12 separate "wonOne".."wonTwelve" boolean columns on the
# "data" sheet with a single numeric "won" column (K), clearing out the
# now-unused L:V columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Header row: K1 becomes "won"; L1:V1 no longer have headers.
$ws.Range("K1").Value = "won"
$ws.Range("L1:V1").ClearContents()

# Row 2 previously had its "won" flag ("y") in column L (wonTwo).
$ws.Range("K2").Value = 7
$ws.Range("L2").ClearContents()

# Row 3 previously had its "won" flag ("y") in column K (wonOne).
$ws.Range("K3").Value = 6

# Row 4 previously had its "won" flag ("y") in column M (wonThree).
# K4 had no prior cell, so copy the shared column formatting (from J4)
# before writing the value.
$ws.Range("J4").Copy()
$ws.Range("K4").PasteSpecial(-4122)
$ws.Range("K4").Value = 12
$ws.Range("M4").ClearContents()

# Row 5 previously had its "won" flag ("y") in column N (wonFour).
$ws.Range("J5").Copy()
$ws.Range("K5").PasteSpecial(-4122)
$ws.Range("K5").Value = 5
$ws.Range("N5").ClearContents()

# Row 6 previously had its "won" flag ("y") in column O (wonFive).
$ws.Range("J6").Copy()
$ws.Range("K6").PasteSpecial(-4122)
$ws.Range("K6").Value = 4
$ws.Range("O6").ClearContents()

# Rows 7-13 previously had their "won" flag ("y") in columns P..V
# (wonSix..wonTwelve). No replacement numeric value was entered yet.
$ws.Range("P7").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S10").ClearContents()
$ws.Range("T11").ClearContents()
$ws.Range("U12").ClearContents()
$ws.Range("V13").ClearContents()
